$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handback"
#
# For each language sheet (zh-cn, de-de) the two data rows (row 2 = the
# 9240c555... file, row 3 = the c1295f7f... file) now report a completed
# handback:
#   - Status (col C)                      -> "Handed back: in sync with en-US"
#   - Latest Target File (col F, new)     -> same source file as col A
#   - Latest Handback File (col G, new)   -> same handoff file as col D
#   - Latest Handback DateTime (col H)    -> an actual timestamp instead of
#                                             the zero-date placeholder
# ---------------------------------------------------------------------------

$statusText    = "Handed back: in sync with en-US"
$handbackStamp = "2016-03-19 10:48:28"

$mdUrl9240 = "https://github.com/OpenLocalizationTest/oltest/blob/ff809d178cca398ed6f917fd52878795e70a97fd/e2e/9240c555-2b70-4f0f-a445-65282756f9e2.md"
$mdUrlc129 = "https://github.com/OpenLocalizationTest/oltest/blob/ff809d178cca398ed6f917fd52878795e70a97fd/e2e/c1295f7f-694d-411f-90ed-e43bd538d728.md"

$mdName9240 = "9240c555-2b70-4f0f-a445-65282756f9e2.md"
$mdNamec129 = "c1295f7f-694d-411f-90ed-e43bd538d728.md"

$rows = @{
    "zh-cn" = @{
        2 = @{
            XlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bbbaadc85401e702cc3b849800550713495392dd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9240c555-2b70-4f0f-a445-65282756f9e2.25b3d082bd92a53a9239577b558616c30f3045c3.zh-cn.xlf"
            XlfName = "9240c555-2b70-4f0f-a445-65282756f9e2.25b3d082bd92a53a9239577b558616c30f3045c3.zh-cn.xlf"
            MdUrl   = $mdUrl9240
            MdName  = $mdName9240
        }
        3 = @{
            XlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bbbaadc85401e702cc3b849800550713495392dd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c1295f7f-694d-411f-90ed-e43bd538d728.406e30fec38b2ef868411948b380b24658f3ed98.zh-cn.xlf"
            XlfName = "c1295f7f-694d-411f-90ed-e43bd538d728.406e30fec38b2ef868411948b380b24658f3ed98.zh-cn.xlf"
            MdUrl   = $mdUrlc129
            MdName  = $mdNamec129
        }
    }
    "de-de" = @{
        2 = @{
            XlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e08f74003158067e878eb1e55de13c2d54753f8a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9240c555-2b70-4f0f-a445-65282756f9e2.25b3d082bd92a53a9239577b558616c30f3045c3.de-de.xlf"
            XlfName = "9240c555-2b70-4f0f-a445-65282756f9e2.25b3d082bd92a53a9239577b558616c30f3045c3.de-de.xlf"
            MdUrl   = $mdUrl9240
            MdName  = $mdName9240
        }
        3 = @{
            XlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e08f74003158067e878eb1e55de13c2d54753f8a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c1295f7f-694d-411f-90ed-e43bd538d728.406e30fec38b2ef868411948b380b24658f3ed98.de-de.xlf"
            XlfName = "c1295f7f-694d-411f-90ed-e43bd538d728.406e30fec38b2ef868411948b380b24658f3ed98.de-de.xlf"
            MdUrl   = $mdUrlc129
            MdName  = $mdNamec129
        }
    }
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in @(2, 3)) {
        $info = $rows[$sheetName][$row]

        # --- Status -------------------------------------------------------
        $ws.Range("C$row").Value = $statusText

        # --- Latest Target File (F) = same file/link as Source File (A) --
        $ws.Hyperlinks.Add($ws.Range("F$row"), $info.MdUrl, "", "", $info.MdName)

        # --- Latest Handback File (G) = same file/link as Handoff File (D)
        $ws.Hyperlinks.Add($ws.Range("G$row"), $info.XlfUrl, "", "", $info.XlfName)

        # --- Latest Handback DateTime (H) ---------------------------------
        $ws.Range("H$row").Value = $handbackStamp
    }
}
